$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.501.04'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '2.110.66'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.10'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5228'
$ws.Range("E7").Value = '  -1.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4499'
$ws.Range("E8").Value = '  +2.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.45'
$ws.Range("E9").Value = '  +15.47%  '
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.164'
$ws.Range("E11").Value = '  -1.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.47'
$ws.Range("E12").Value = '  -2.60%  '
$ws.Range("D13").Value = '2.103.10'
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.790'
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.834'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.58'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06610'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.306'
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").Value = '30.548.42'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.36'
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("E25").Value = '  +3.38%  '
$ws.Range("D26").Value = '2.348.18'
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.586'
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.90'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.97'
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.203'
$ws.Range("E31").Value = '  +2.41%  '
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.670'
$ws.Range("E33").Value = '  +7.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.170'
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.942'
$ws.Range("E35").Value = '  -1.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.58'
$ws.Range("E36").Value = '  +11.28%  '
$ws.Range("E37").Value = '  -1.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06797'
$ws.Range("E38").Value = '  +0.69%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.510'
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.76'
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2279'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6934'
$ws.Range("E42").Value = '  +1.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.254'
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.351'
$ws.Range("E44").Value = '  +5.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '14.09'
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6395'
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.652'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.246'
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.225'
$ws.Range("E50").Value = '  +5.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '83.36'
$ws.Range("E51").Value = '  +0.39%  '
